$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "fixed motor efficiency glitch" -----------------------------------
# The FC list had AliceBallard (Cell Number 46) in row 2 and HalfBallard
# (Cell Number 23) in row 3. The fix moves HalfBallard up into row 2 and
# pushes AliceBallard down to row 8 (leaving row 7 blank), making room for
# four newly added Ballard fuel cells in rows 3-6.

# Seed the AliceBallard text in its new home (row 8) first, before row 2's
# text is overwritten, so the "AliceBallard" shared string is preserved.
$ws.Range("A8").Value = "AliceBallard"
$ws.Range("B8").Value = 46
$ws.Range("C8").Value = 145
$ws.Range("D8").Value = 0.36
$ws.Range("E8").Value = 0.45
$ws.Range("F8").Value = 0.04
$ws.Range("G8").Value = 1.02

# Row 2: becomes HalfBallard's row (was row 3's data).
$ws.Range("A2").Value = "HalfBallard"
$ws.Range("B2").Value = 23

# Row 3: becomes the new Ballard22 fuel cell.
$ws.Range("A3").Value = "Ballard22"
$ws.Range("B3").Value = 22

# --- "added more ... fuelcells" ----------------------------------------
# New fuel cell rows 4-6: Ballard21, Ballard24, Ballard25.
$ws.Range("A4").Value = "Ballard21"
$ws.Range("B4").Value = 21
$ws.Range("C4").Value = 145
$ws.Range("D4").Value = 0.36
$ws.Range("E4").Value = 0.45
$ws.Range("F4").Value = 0.04
$ws.Range("G4").Value = 1.02

$ws.Range("A5").Value = "Ballard24"
$ws.Range("B5").Value = 24
$ws.Range("C5").Value = 145
$ws.Range("D5").Value = 0.36
$ws.Range("E5").Value = 0.45
$ws.Range("F5").Value = 0.04
$ws.Range("G5").Value = 1.02

$ws.Range("A6").Value = "Ballard25"
$ws.Range("B6").Value = 25
$ws.Range("C6").Value = 145
$ws.Range("D6").Value = 0.36
$ws.Range("E6").Value = 0.45
$ws.Range("F6").Value = 0.04
$ws.Range("G6").Value = 1.02

# Row 7 is intentionally left blank.
$ws.Range("A7").Select()
